# add fixedtop settings for masthead
#
# For the OVERVIEW sheet, append two new rows (icon/fa-star, slug/index)
# after the existing 4 rows (topic/headline/byline/reporter_bio).
#
# For every topic sheet (HEALTH, PUBLIC_EDU, HIGHER_EDU, TRANSPO,
# IMMIGRATION, ENERGY, ENVIRO, TEF, JUSTICE) append a new row 6 with
# slug/<topic-slug> after the existing 5 rows (...,icon/fa-*).

$wb = $excel.ActiveWorkbook

# NOTE: this interpreter only reliably binds POSITIONAL function arguments
# (named `-Param value` binding was observed to pass $null through), so
# Add-KeyValueRow is called positionally everywhere below.
function Add-KeyValueRow {
    param($Worksheet, $RowNumber, $Key, $Value)

    # Copy formatting from row 1 (plain, non-hyperlink style) down onto the
    # new row so the new cells pick up the same cell style ("s") as the
    # rest of column A / B, then overwrite with the real key/value text.
    $srcRange = $Worksheet.Range("A1:B1")
    $srcRange.Copy() | Out-Null

    $dstRange = $Worksheet.Range("A" + $RowNumber + ":B" + $RowNumber)
    $dstRange.PasteSpecial(-4122) | Out-Null

    $Worksheet.Range("A" + $RowNumber).Value = $Key
    $Worksheet.Range("B" + $RowNumber).Value = $Value
}

# OVERVIEW: add icon/fa-star (row 5) and slug/index (row 6)
$wsOverview = $wb.Worksheets.Item("OVERVIEW")
Add-KeyValueRow $wsOverview 5 "icon" "fa-star"
Add-KeyValueRow $wsOverview 6 "slug" "index"

# Topic sheets: each gets a new row 6 with slug/<slug-value>
$slugBySheet = [ordered]@{
    "HEALTH"      = "health-care"
    "PUBLIC_EDU"  = "public-education"
    "HIGHER_EDU"  = "higher-education"
    "TRANSPO"     = "transportation"
    "IMMIGRATION" = "immigration"
    "ENERGY"      = "energy"
    "ENVIRO"      = "environment"
    "TEF"         = "texas-enterprise-fund"
    "JUSTICE"     = "criminal-justice"
}

foreach ($sheetName in $slugBySheet.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    Add-KeyValueRow $ws 6 "slug" $slugBySheet[$sheetName]
}
